$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.4841736666666667
$ws.Range("H2").Value = 1.452521
$ws.Range("I2").Value = 0.003342651198412304
$ws.Range("J2").Value = 0.003342651198412304
$ws.Range("M2").Value = 2.808016
$ws.Range("N2").Value = 8.424047999999999
$ws.Range("O2").Value = 0.02933065711877975
$ws.Range("P2").Value = 0.02933065711877976
$ws.Range("Q2").Value = 1.359567402778667
$ws.Range("R2").Value = 12.236106625008
$ws.Range("S2").Value = 0.00009804215616830952
$ws.Range("T2").Value = 0.00009804215616830953
$ws.Range("G3").Value = 0.4841736666666667
$ws.Range("H3").Value = 1.452521
$ws.Range("I3").Value = 0.003342651198412304
$ws.Range("J3").Value = 0.003342651198412304
$ws.Range("O3").Value = 0.05404782035567898
$ws.Range("P3").Value = 0.05404782035567899
$ws.Range("Q3").Value = 2.505284980464
$ws.Range("R3").Value = 22.547564824176
$ws.Range("S3").Value = 0.0001806630114834833
$ws.Range("T3").Value = 0.0001806630114834833
$ws.Range("G4").Value = 0.4841736666666667
$ws.Range("H4").Value = 1.452521
$ws.Range("I4").Value = 0.003342651198412304
$ws.Range("J4").Value = 0.003342651198412304
$ws.Range("M4").Value = 50.87875466666667
$ws.Range("N4").Value = 152.636264
$ws.Range("O4").Value = 0.5314454432448089
$ws.Range("P4").Value = 0.531445443244809
$ws.Range("Q4").Value = 24.63415320239378
$ws.Range("R4").Value = 221.707378821544
$ws.Range("S4").Value = 0.001776436747753019
$ws.Range("T4").Value = 0.001776436747753019
$ws.Range("G5").Value = 0.4841736666666667
$ws.Range("H5").Value = 1.452521
$ws.Range("I5").Value = 0.003342651198412304
$ws.Range("J5").Value = 0.003342651198412304
$ws.Range("M5").Value = 0.7109426666666666
$ws.Range("N5").Value = 2.132828
$ws.Range("O5").Value = 0.007426031613463359
$ws.Range("P5").Value = 0.00742603161346336
$ws.Range("Q5").Value = 0.3442197177097778
$ws.Range("R5").Value = 3.097977459388
$ws.Range("S5").Value = 0.00002482263347219095
$ws.Range("T5").Value = 0.00002482263347219096
$ws.Range("G6").Value = 0.4841736666666667
$ws.Range("H6").Value = 1.452521
$ws.Range("I6").Value = 0.003342651198412304
$ws.Range("J6").Value = 0.003342651198412304
$ws.Range("M6").Value = 36.16448733333333
$ws.Range("N6").Value = 108.493462
$ws.Range("O6").Value = 0.3777500476672688
$ws.Range("P6").Value = 0.3777500476672689
$ws.Range("Q6").Value = 17.50989243530022
$ws.Range("R6").Value = 157.589031917702
$ws.Range("S6").Value = 0.001262686649535301
$ws.Range("T6").Value = 0.001262686649535301
$ws.Range("I7").Value = 0.00775849308788848
$ws.Range("J7").Value = 0.00775849308788848
$ws.Range("M7").Value = 2.808016
$ws.Range("N7").Value = 8.424047999999999
$ws.Range("O7").Value = 0.02933065711877975
$ws.Range("P7").Value = 0.02933065711877976
$ws.Range("Q7").Value = 3.155637148736
$ws.Range("R7").Value = 28.400734338624
$ws.Range("S7").Value = 0.0002275617005192798
$ws.Range("T7").Value = 0.0002275617005192798
$ws.Range("I8").Value = 0.00775849308788848
$ws.Range("J8").Value = 0.00775849308788848
$ws.Range("O8").Value = 0.05404782035567898
$ws.Range("P8").Value = 0.05404782035567899
$ws.Range("S8").Value = 0.0004193296406449736
$ws.Range("T8").Value = 0.0004193296406449737
$ws.Range("I9").Value = 0.00775849308788848
$ws.Range("J9").Value = 0.00775849308788848
$ws.Range("M9").Value = 50.87875466666667
$ws.Range("N9").Value = 152.636264
$ws.Range("O9").Value = 0.5314454432448089
$ws.Range("P9").Value = 0.531445443244809
$ws.Range("Q9").Value = 57.17734097938136
$ws.Range("R9").Value = 514.5960688144321
$ws.Range("S9").Value = 0.004123215798004679
$ws.Range("T9").Value = 0.00412321579800468
$ws.Range("I10").Value = 0.00775849308788848
$ws.Range("J10").Value = 0.00775849308788848
$ws.Range("M10").Value = 0.7109426666666666
$ws.Range("N10").Value = 2.132828
$ws.Range("O10").Value = 0.007426031613463359
$ws.Range("P10").Value = 0.00742603161346336
$ws.Range("Q10").Value = 0.7989545250293334
$ws.Range("R10").Value = 7.190590725264001
$ws.Range("S10").Value = 0.00005761481494349681
$ws.Range("T10").Value = 0.00005761481494349682
$ws.Range("I11").Value = 0.00775849308788848
$ws.Range("J11").Value = 0.00775849308788848
$ws.Range("M11").Value = 36.16448733333333
$ws.Range("N11").Value = 108.493462
$ws.Range("O11").Value = 0.3777500476672688
$ws.Range("P11").Value = 0.3777500476672689
$ws.Range("Q11").Value = 40.64150620725068
$ws.Range("R11").Value = 365.773555865256
$ws.Range("S11").Value = 0.002930771133776049
$ws.Range("T11").Value = 0.00293077113377605
$ws.Range("G12").Value = 85.95243833333332
$ws.Range("H12").Value = 257.857315
$ws.Range("I12").Value = 0.5934007584084009
$ws.Range("J12").Value = 0.5934007584084009
$ws.Range("M12").Value = 2.808016
$ws.Range("N12").Value = 8.424047999999999
$ws.Range("O12").Value = 0.02933065711877975
$ws.Range("P12").Value = 0.02933065711877976
$ws.Range("Q12").Value = 241.3558220790133
$ws.Range("R12").Value = 2172.20239871112
$ws.Range("S12").Value = 0.01740483417890067
$ws.Range("T12").Value = 0.01740483417890067
$ws.Range("G13").Value = 85.95243833333332
$ws.Range("H13").Value = 257.857315
$ws.Range("I13").Value = 0.5934007584084009
$ws.Range("J13").Value = 0.5934007584084009
$ws.Range("O13").Value = 0.05404782035567898
$ws.Range("P13").Value = 0.05404782035567899
$ws.Range("Q13").Value = 444.7481711949599
$ws.Range("R13").Value = 4002.73354075464
$ws.Range("S13").Value = 0.03207201758938091
$ws.Range("T13").Value = 0.03207201758938092
$ws.Range("G14").Value = 85.95243833333332
$ws.Range("H14").Value = 257.857315
$ws.Range("I14").Value = 0.5934007584084009
$ws.Range("J14").Value = 0.5934007584084009
$ws.Range("M14").Value = 50.87875466666667
$ws.Range("N14").Value = 152.636264
$ws.Range("O14").Value = 0.5314454432448089
$ws.Range("P14").Value = 0.531445443244809
$ws.Range("Q14").Value = 4373.153022963462
$ws.Range("R14").Value = 39358.37720667116
$ws.Range("S14").Value = 0.3153601290741584
$ws.Range("T14").Value = 0.3153601290741584
$ws.Range("G15").Value = 85.95243833333332
$ws.Range("H15").Value = 257.857315
$ws.Range("I15").Value = 0.5934007584084009
$ws.Range("J15").Value = 0.5934007584084009
$ws.Range("M15").Value = 0.7109426666666666
$ws.Range("N15").Value = 2.132828
$ws.Range("O15").Value = 0.007426031613463359
$ws.Range("P15").Value = 0.00742603161346336
$ws.Range("Q15").Value = 61.10725571520221
$ws.Range("R15").Value = 549.9653014368199
$ws.Range("S15").Value = 0.004406612791393918
$ws.Range("T15").Value = 0.004406612791393919
$ws.Range("G16").Value = 85.95243833333332
$ws.Range("H16").Value = 257.857315
$ws.Range("I16").Value = 0.5934007584084009
$ws.Range("J16").Value = 0.5934007584084009
$ws.Range("M16").Value = 36.16448733333333
$ws.Range("N16").Value = 108.493462
$ws.Range("O16").Value = 0.3777500476672688
$ws.Range("P16").Value = 0.3777500476672689
$ws.Range("Q16").Value = 3108.425867374947
$ws.Range("R16").Value = 27975.83280637452
$ws.Range("S16").Value = 0.2241571647745669
$ws.Range("T16").Value = 0.224157164774567
$ws.Range("G17").Value = 1.941884333333334
$ws.Range("H17").Value = 5.825653000000001
$ws.Range("I17").Value = 0.01340643335413687
$ws.Range("J17").Value = 0.01340643335413687
$ws.Range("M17").Value = 2.808016
$ws.Range("N17").Value = 8.424047999999999
$ws.Range("O17").Value = 0.02933065711877975
$ws.Range("P17").Value = 0.02933065711877976
$ws.Range("Q17").Value = 5.452842278149333
$ws.Range("R17").Value = 49.075580503344
$ws.Range("S17").Value = 0.0003932194998959608
$ws.Range("T17").Value = 0.0003932194998959609
$ws.Range("G18").Value = 1.941884333333334
$ws.Range("H18").Value = 5.825653000000001
$ws.Range("I18").Value = 0.01340643335413687
$ws.Range("J18").Value = 0.01340643335413687
$ws.Range("O18").Value = 0.05404782035567898
$ws.Range("P18").Value = 0.05404782035567899
$ws.Range("Q18").Value = 10.047993083952
$ws.Range("R18").Value = 90.43193775556801
$ws.Range("S18").Value = 0.0007245885015347722
$ws.Range("T18").Value = 0.0007245885015347723
$ws.Range("G19").Value = 1.941884333333334
$ws.Range("H19").Value = 5.825653000000001
$ws.Range("I19").Value = 0.01340643335413687
$ws.Range("J19").Value = 0.01340643335413687
$ws.Range("M19").Value = 50.87875466666667
$ws.Range("N19").Value = 152.636264
$ws.Range("O19").Value = 0.5314454432448089
$ws.Range("P19").Value = 0.531445443244809
$ws.Range("Q19").Value = 98.80065658671025
$ws.Range("R19").Value = 889.2059092803922
$ws.Range("S19").Value = 0.007124787916221258
$ws.Range("T19").Value = 0.00712478791622126
$ws.Range("G20").Value = 1.941884333333334
$ws.Range("H20").Value = 5.825653000000001
$ws.Range("I20").Value = 0.01340643335413687
$ws.Range("J20").Value = 0.01340643335413687
$ws.Range("M20").Value = 0.7109426666666666
$ws.Range("N20").Value = 2.132828
$ws.Range("O20").Value = 0.007426031613463359
$ws.Range("P20").Value = 0.00742603161346336
$ws.Range("Q20").Value = 1.380568426298222
$ws.Range("R20").Value = 12.425115836684
$ws.Range("S20").Value = 0.00009955659791160999
$ws.Range("T20").Value = 0.00009955659791161001
$ws.Range("G21").Value = 1.941884333333334
$ws.Range("H21").Value = 5.825653000000001
$ws.Range("I21").Value = 0.01340643335413687
$ws.Range("J21").Value = 0.01340643335413687
$ws.Range("M21").Value = 36.16448733333333
$ws.Range("N21").Value = 108.493462
$ws.Range("O21").Value = 0.3777500476672688
$ws.Range("P21").Value = 0.3777500476672689
$ws.Range("Q21").Value = 70.22725137563178
$ws.Range("R21").Value = 632.045262380686
$ws.Range("S21").Value = 0.005064280838573265
$ws.Range("T21").Value = 0.005064280838573266
$ws.Range("G22").Value = 55.34490766666666
$ws.Range("H22").Value = 166.034723
$ws.Range("I22").Value = 0.3820916639511615
$ws.Range("J22").Value = 0.3820916639511614
$ws.Range("M22").Value = 2.808016
$ws.Range("N22").Value = 8.424047999999999
$ws.Range("O22").Value = 0.02933065711877975
$ws.Range("P22").Value = 0.02933065711877976
$ws.Range("Q22").Value = 155.4093862465226
$ws.Range("R22").Value = 1398.684476218704
$ws.Range("S22").Value = 0.01120699958329554
$ws.Range("T22").Value = 0.01120699958329553
$ws.Range("G23").Value = 55.34490766666666
$ws.Range("H23").Value = 166.034723
$ws.Range("I23").Value = 0.3820916639511615
$ws.Range("J23").Value = 0.3820916639511614
$ws.Range("O23").Value = 0.05404782035567898
$ws.Range("P23").Value = 0.05404782035567899
$ws.Range("Q23").Value = 286.374033674832
$ws.Range("R23").Value = 2577.366303073488
$ws.Range("S23").Value = 0.02065122161263484
$ws.Range("T23").Value = 0.02065122161263484
$ws.Range("G24").Value = 55.34490766666666
$ws.Range("H24").Value = 166.034723
$ws.Range("I24").Value = 0.3820916639511615
$ws.Range("J24").Value = 0.3820916639511614
$ws.Range("M24").Value = 50.87875466666667
$ws.Range("N24").Value = 152.636264
$ws.Range("O24").Value = 0.5314454432448089
$ws.Range("P24").Value = 0.531445443244809
$ws.Range("Q24").Value = 2815.879979221653
$ws.Range("R24").Value = 25342.91981299487
$ws.Range("S24").Value = 0.2030608737086716
$ws.Range("T24").Value = 0.2030608737086716
$ws.Range("G25").Value = 55.34490766666666
$ws.Range("H25").Value = 166.034723
$ws.Range("I25").Value = 0.3820916639511615
$ws.Range("J25").Value = 0.3820916639511614
$ws.Range("M25").Value = 0.7109426666666666
$ws.Range("N25").Value = 2.132828
$ws.Range("O25").Value = 0.007426031613463359
$ws.Range("P25").Value = 0.00742603161346336
$ws.Range("Q25").Value = 39.34705624296044
$ws.Range("R25").Value = 354.123506186644
$ws.Range("S25").Value = 0.002837424775742143
$ws.Range("T25").Value = 0.002837424775742143
$ws.Range("G26").Value = 55.34490766666666
$ws.Range("H26").Value = 166.034723
$ws.Range("I26").Value = 0.3820916639511615
$ws.Range("J26").Value = 0.3820916639511614
$ws.Range("M26").Value = 36.16448733333333
$ws.Range("N26").Value = 108.493462
$ws.Range("O26").Value = 0.3777500476672688
$ws.Range("P26").Value = 0.3777500476672689
$ws.Range("Q26").Value = 2001.52021227567
$ws.Range("R26").Value = 18013.68191048102
$ws.Range("S26").Value = 0.1443351442708173
$ws.Range("T26").Value = 0.1443351442708173
